$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly record per row (rows 2..122). A new weekly
# record is being inserted as the new "row 51" (date 2023-10-04), which
# pushes every existing record from old-row 51 down through old-row 122
# down by one row (new rows 52..123). Columns A,B,C,E,F,G,H,N,Q,R are
# constant for every data row, so only D (Fecha), I (Calidad), J
# (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado), O (Origen) and P (Precio $/Kg) actually need to move.

# Row 123 lands outside the sheet's previous used range, so its Fecha
# cell has no number format yet - give it the same date format as every
# other cell in column D before writing the date value into it (so Excel
# doesn't invent a fresh "m/d/yyyy" style for it).
$ws.Range("D123").NumberFormat = $ws.Range("D122").NumberFormat()

# 1) Shift every data row from old-51..old-122 down to new-52..new-123 (a
# full-row copy is simplest and correct since columns A,B,C,E,F,G,H,N,Q,R
# are identical for every record anyway).
$ws.Range("A52:R123").Value = $ws.Range("A51:R122").Value()

# 2) Overwrite row 51 with the new weekly record. Calidad (I51) and
# Volumen (J51) are unchanged; Fecha, the three prices and Origen change.
$ws.Range("D51").Value = "2023-10-04"
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = 2000
$ws.Range("O51").Value = "Región de Ñuble"
$ws.Range("P51").Value = 2000
